$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.787.53"
$ws.Range("E2").Value = "  +7.99%  "
$ws.Range("D3").Value = "1.953.16"
$ws.Range("E3").Value = "  +6.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "341.53"
$ws.Range("E5").Value = "  +2.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4777"
$ws.Range("E7").Value = "  +3.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4144"
$ws.Range("E8").Value = "  +7.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.86"
$ws.Range("E9").Value = "  +3.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08244"
$ws.Range("E10").Value = "  +4.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.037"
$ws.Range("E11").Value = "  +7.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.77"
$ws.Range("E12").Value = "  +7.70%  "
$ws.Range("D13").Value = "1.952.09"
$ws.Range("E13").Value = "  +5.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.171"
$ws.Range("E14").Value = "  +5.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.398"
$ws.Range("E15").Value = "  +4.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.03"
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("E18").Value = "  +3.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06698"
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.04"
$ws.Range("E20").Value = "  +5.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "29.776.77"
$ws.Range("E22").Value = "  +7.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.600"
$ws.Range("E23").Value = "  +5.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.27"
$ws.Range("E24").Value = "  +4.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.278"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("D26").Value = "2.179.59"
$ws.Range("E26").Value = "  +5.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.18"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.28"
$ws.Range("E28").Value = "  +4.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.186"
$ws.Range("E29").Value = "  +6.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.699"
$ws.Range("E30").Value = "  +7.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.14"
$ws.Range("E31").Value = "  +4.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.004"
$ws.Range("E32").Value = "  +7.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09650"
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.476"
$ws.Range("E34").Value = "  +11.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.685"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.518"
$ws.Range("E36").Value = "  +5.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06269"
$ws.Range("E37").Value = "  +5.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02315"
$ws.Range("E38").Value = "  +5.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.439"
$ws.Range("E39").Value = "  +3.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6104"
$ws.Range("E40").Value = "  +5.89%  "
$ws.Range("E41").Value = "  +3.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.74"
$ws.Range("E42").Value = "  +7.70%  "
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1893"
$ws.Range("E44").Value = "  +3.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.279"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.384"
$ws.Range("E46").Value = "  +33.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5716"
$ws.Range("E47").Value = "  +5.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.50"
$ws.Range("E48").Value = "  +5.95%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.992"
$ws.Range("E49").Value = "  +4.68%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07392"
$ws.Range("E50").Value = "  +8.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.13"
$ws.Range("E51").Value = "  +1.90%  "
